# Update countries & provincias Spain
# - Swap display order of "Libia" / "San Martin (Parte Holandesa)" rows
#   (country table is sorted by case count, so the two countries swap
#   table rows as their numbers change).
# - Refresh the "Datos actualizados..." timestamp string.
# - Refresh the per-country case counters for the rows whose figures moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp banner (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 20:05"

# --- country name swap: row 177 <-> row 178 --------------------------------
$ws.Range("A177").Value = "Libia"
$ws.Range("A178").Value = "San Martin (Parte Holandesa)"

# --- row 4 (rank 8) ---------------------------------------------------------
$ws.Range("B4").Value = 1714327
$ws.Range("C4").Value = 8101
$ws.Range("D4").Value = 469049
$ws.Range("E4").Value = 1145175
$ws.Range("G4").Value = 298
$ws.Range("H4").Value = 100103

# --- row 5 (rank 9) ---------------------------------------------------------
$ws.Range("B5").Value = 377780
$ws.Range("C5").Value = 1111
$ws.Range("E5").Value = 200325
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 23622

# --- row 13 (rank 17) --------------------------------------------------------
$ws.Range("B13").Value = 150762
$ws.Range("C13").Value = 5812
$ws.Range("E13").Value = 82141

# --- row 53 (rank 57) --------------------------------------------------------
$ws.Range("B53").Value = 9366
$ws.Range("C53").Value = 195
$ws.Range("D53").Value = 4938
$ws.Range("E53").Value = 4414

# --- row 177 (rank 181, now "Libia") -----------------------------------------
$ws.Range("C177").Value = 2
$ws.Range("D177").Value = 40
$ws.Range("E177").Value = 34
$ws.Range("H177").Value = 3

# --- row 178 (rank 182, now "San Martin (Parte Holandesa)") ------------------
$ws.Range("B178").Value = 77
$ws.Range("D178").Value = 59
$ws.Range("E178").Value = 3
$ws.Range("H178").Value = 15
